$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Write discount value
$ws.Range("E18").Value = 100

# Write Signature (name, email address)
$ws.Range("A31").Value = "RPA Dev, developer.rpa@mail.com"
